$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.819.13"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "3.119.62"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.60"
$ws.Range("E5").Value = "  +2.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.63"
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.492"
$ws.Range("E8").Value = "  +9.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.35"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.412"
$ws.Range("E11").Value = "  +3.61%  "
$ws.Range("E12").Value = "  +3.88%  "
$ws.Range("D13").Value = "3.658.67"
$ws.Range("E13").Value = "  +1.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.90"
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000166"
$ws.Range("E15").Value = "  +3.15%  "
$ws.Range("D16").Value = "57.945.88"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("D17").Value = "3.121.08"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("E18").Value = "  +4.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.81"
$ws.Range("E19").Value = "  +2.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.12"
$ws.Range("E20").Value = "  +3.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "372.86"
$ws.Range("E21").Value = "  +7.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("E23").Value = "  -1.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.31"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.508"
$ws.Range("E25").Value = "  +2.24%  "
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").Value = "0.0₃0880"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  +4.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.17"
$ws.Range("E30").Value = "  +4.57%  "
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.54"
$ws.Range("E32").Value = "  +3.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.15"
$ws.Range("E33").Value = "  +4.85%  "
$ws.Range("E34").Value = "  +3.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "160.29"
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.17"
$ws.Range("E36").Value = "  +2.79%  "
$ws.Range("E37").Value = "  +4.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.59"
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.64"
$ws.Range("E39").Value = "  +3.93%  "
$ws.Range("E40").Value = "  +2.87%  "
$ws.Range("D41").Value = "2.561.13"
$ws.Range("E42").Value = "  +3.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.23"
$ws.Range("E43").Value = "  +4.53%  "
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("E45").Value = "  +3.17%  "
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.981"
$ws.Range("E47").Value = "  +2.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.15"
$ws.Range("E48").Value = "  +3.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.88"
$ws.Range("E49").Value = "  +1.34%  "
$ws.Range("E50").Value = "  +5.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.747"
$ws.Range("E51").Value = "  -0.88%  "
